# Scheduled runner update: refresh cached Universalis market-price data
# (currentAveragePrice / LevePrice / LeveProfit columns) for affected leve rows
# across the per-job Pandaemonium profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 325
$ws.Range("I18").Value = 315.27777
$ws.Range("K18").Value = 315.27777
$ws.Range("M18").Value = -31.27776999999998

$ws.Range("H74").Value = 4836.9165
$ws.Range("I74").Value = 4504.3
$ws.Range("J74").Value = 6500
$ws.Range("K74").Value = 4504.3
$ws.Range("L74").Value = 6500
$ws.Range("M74").Value = -3568.3
$ws.Range("N74").Value = -8372

$ws.Range("H77").Value = 4836.9165
$ws.Range("I77").Value = 4504.3
$ws.Range("J77").Value = 6500
$ws.Range("K77").Value = 22521.5
$ws.Range("L77").Value = 32500
$ws.Range("M77").Value = -17841.5
$ws.Range("N77").Value = -41860

$ws.Range("H112").Value = 6582.2856
$ws.Range("I112").Value = 70100
$ws.Range("J112").Value = 1696.3077
$ws.Range("K112").Value = 210300
$ws.Range("L112").Value = 5088.9231
$ws.Range("M112").Value = -209192
$ws.Range("N112").Value = -7304.9231

$ws.Range("H137").Value = 816093.9
$ws.Range("I137").Value = 3688.5
$ws.Range("J137").Value = 1336033.2
$ws.Range("K137").Value = 11065.5
$ws.Range("L137").Value = 4008099.6
$ws.Range("M137").Value = -8515.5
$ws.Range("N137").Value = -4013199.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16234.395
$ws.Range("I32").Value = 18260.951
$ws.Range("J32").Value = 3872.4
$ws.Range("K32").Value = 18260.951
$ws.Range("L32").Value = 3872.4
$ws.Range("M32").Value = -17973.951
$ws.Range("N32").Value = -4446.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 399.6
$ws.Range("I64").Value = 374.5
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 374.5
$ws.Range("L64").Value = 500
$ws.Range("M64").Value = -149.5
$ws.Range("N64").Value = -950

$ws.Range("H67").Value = 399.6
$ws.Range("I67").Value = 374.5
$ws.Range("J67").Value = 500
$ws.Range("K67").Value = 374.5
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 405.5
$ws.Range("N67").Value = -2060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 323.33334
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000

$ws.Range("H31").Value = 591056.4
$ws.Range("I31").Value = 5219.913
$ws.Range("J31").Value = 927912.3
$ws.Range("K31").Value = 5219.913
$ws.Range("L31").Value = 927912.3
$ws.Range("M31").Value = -4924.913
$ws.Range("N31").Value = -928502.3

$ws.Range("H34").Value = 591056.4
$ws.Range("I34").Value = 5219.913
$ws.Range("J34").Value = 927912.3
$ws.Range("K34").Value = 5219.913
$ws.Range("L34").Value = 927912.3
$ws.Range("M34").Value = -5017.913
$ws.Range("N34").Value = -928316.3

$ws.Range("H41").Value = 31954.666
$ws.Range("J41").Value = 17932.5
$ws.Range("L41").Value = 17932.5
$ws.Range("N41").Value = -18788.5

$ws.Range("H50").Value = 26652.637
$ws.Range("J50").Value = 26652.637
$ws.Range("L50").Value = 26652.637
$ws.Range("N50").Value = -27902.637

$ws.Range("H51").Value = 20783.8
$ws.Range("J51").Value = 20783.8
$ws.Range("L51").Value = 20783.8
$ws.Range("N51").Value = -22255.8

$ws.Range("H59").Value = 25115.5
$ws.Range("J59").Value = 40127
$ws.Range("L59").Value = 40127
$ws.Range("N59").Value = -42417

$ws.Range("H60").Value = 20223.25
$ws.Range("J60").Value = 23600
$ws.Range("L60").Value = 23600
$ws.Range("N60").Value = -24622

$ws.Range("H61").Value = 20783.8
$ws.Range("J61").Value = 20783.8
$ws.Range("L61").Value = 20783.8
$ws.Range("N61").Value = -21479.8

$ws.Range("H62").Value = 3975.125
$ws.Range("I62").Value = 3685
$ws.Range("J62").Value = 6006
$ws.Range("K62").Value = 3685
$ws.Range("L62").Value = 6006
$ws.Range("M62").Value = -3061
$ws.Range("N62").Value = -7254

$ws.Range("H65").Value = 3975.125
$ws.Range("I65").Value = 3685
$ws.Range("J65").Value = 6006
$ws.Range("K65").Value = 18425
$ws.Range("L65").Value = 30030
$ws.Range("M65").Value = -15305
$ws.Range("N65").Value = -36270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3010.3467
$ws.Range("I68").Value = 1560.9744
$ws.Range("J68").Value = 4580.5
$ws.Range("K68").Value = 4682.9232
$ws.Range("L68").Value = 13741.5
$ws.Range("M68").Value = -3871.9232
$ws.Range("N68").Value = -15363.5

$ws.Range("H71").Value = 3010.3467
$ws.Range("I71").Value = 1560.9744
$ws.Range("J71").Value = 4580.5
$ws.Range("K71").Value = 14048.7696
$ws.Range("L71").Value = 41224.5
$ws.Range("M71").Value = -9992.7696
$ws.Range("N71").Value = -49336.5

$ws.Range("H107").Value = 586.6163
$ws.Range("J107").Value = 1249.3334
$ws.Range("L107").Value = 3748.0002
$ws.Range("N107").Value = -7588.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 10000
$ws.Range("J53").Value = 10000
$ws.Range("L53").Value = 10000
$ws.Range("N53").Value = -11262

$ws.Range("H80").Value = 6086.9375
$ws.Range("I80").Value = 7976.6665
$ws.Range("J80").Value = 3657.2856
$ws.Range("K80").Value = 7976.6665
$ws.Range("L80").Value = 3657.2856
$ws.Range("M80").Value = -6978.6665
$ws.Range("N80").Value = -5653.2856

$ws.Range("H83").Value = 6086.9375
$ws.Range("I83").Value = 7976.6665
$ws.Range("J83").Value = 3657.2856
$ws.Range("K83").Value = 39883.3325
$ws.Range("L83").Value = 18286.428
$ws.Range("M83").Value = -34891.3325
$ws.Range("N83").Value = -28270.428

$ws.Range("H113").Value = 1932.9565
$ws.Range("J113").Value = 2075.5557
$ws.Range("L113").Value = 2075.5557
$ws.Range("N113").Value = -6415.5557

$ws.Range("H126").Value = 2233.2856
$ws.Range("I126").Value = 1676.65
$ws.Range("J126").Value = 2975.4666
$ws.Range("K126").Value = 5029.950000000001
$ws.Range("L126").Value = 8926.399800000001
$ws.Range("M126").Value = -2559.950000000001
$ws.Range("N126").Value = -13866.3998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 140.25
$ws.Range("I55").Value = 137
$ws.Range("J55").Value = 150
$ws.Range("K55").Value = 137
$ws.Range("L55").Value = 150
$ws.Range("M55").Value = 36
$ws.Range("N55").Value = -496

$ws.Range("H61").Value = 25567.549
$ws.Range("I61").Value = 26068.96
$ws.Range("J61").Value = 22960.2
$ws.Range("K61").Value = 26068.96
$ws.Range("L61").Value = 22960.2
$ws.Range("M61").Value = -25866.96
$ws.Range("N61").Value = -23364.2

$ws.Range("H68").Value = 3399.8572
$ws.Range("I68").Value = 2399.8
$ws.Range("J68").Value = 5900
$ws.Range("K68").Value = 2399.8
$ws.Range("L68").Value = 5900
$ws.Range("M68").Value = -1650.8
$ws.Range("N68").Value = -7398

$ws.Range("H71").Value = 3399.8572
$ws.Range("I71").Value = 2399.8
$ws.Range("J71").Value = 5900
$ws.Range("K71").Value = 11999
$ws.Range("L71").Value = 29500
$ws.Range("M71").Value = -8255
$ws.Range("N71").Value = -36988

$ws.Range("H113").Value = 25567.549
$ws.Range("I113").Value = 26068.96
$ws.Range("J113").Value = 22960.2
$ws.Range("K113").Value = 26068.96
$ws.Range("L113").Value = 22960.2
$ws.Range("M113").Value = -23898.96
$ws.Range("N113").Value = -27300.2
